$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "Add 1 to a number represented as linked list"
$ws.Range("A4").Value = "GFG"

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null

$ws.Range("B9").Select() | Out-Null
